$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: repeat the category header labels from row 6 (adm, police, education, health, social, military)
$ws.Cells.Item(34, 5).Value = "adm"
$ws.Cells.Item(34, 6).Value = "police"
$ws.Cells.Item(34, 7).Value = "education"
$ws.Cells.Item(34, 8).Value = "health"
$ws.Cells.Item(34, 9).Value = "social"
$ws.Cells.Item(34, 10).Value = "military"

# Rows 35-58: expected spending data per ideology
$ws.Cells.Item(35, 4).Value = "#Western_Autocracy"
$ws.Cells.Item(35, 5).Value = 2.5
$ws.Cells.Item(35, 6).Value = 3.5
$ws.Cells.Item(35, 7).Value = 2
$ws.Cells.Item(35, 8).Value = 2
$ws.Cells.Item(35, 9).Value = 2
$ws.Cells.Item(35, 10).Value = 2

$ws.Cells.Item(36, 4).Value = "#conservatism"
$ws.Cells.Item(36, 5).Value = 2.5
$ws.Cells.Item(36, 6).Value = 3.5
$ws.Cells.Item(36, 7).Value = 2
$ws.Cells.Item(36, 8).Value = 1.5
$ws.Cells.Item(36, 9).Value = 2
$ws.Cells.Item(36, 10).Value = 2.5

$ws.Cells.Item(37, 4).Value = "#liberalism"
$ws.Cells.Item(37, 5).Value = 1.5
$ws.Cells.Item(37, 6).Value = 2
$ws.Cells.Item(37, 7).Value = 4
$ws.Cells.Item(37, 8).Value = 2.5
$ws.Cells.Item(37, 9).Value = 2
$ws.Cells.Item(37, 10).Value = 2

$ws.Cells.Item(38, 4).Value = "#socialism"
$ws.Cells.Item(38, 5).Value = 2
$ws.Cells.Item(38, 6).Value = 1
$ws.Cells.Item(38, 7).Value = 3
$ws.Cells.Item(38, 8).Value = 3
$ws.Cells.Item(38, 9).Value = 3
$ws.Cells.Item(38, 10).Value = 2

$ws.Cells.Item(39, 4).Value = "Communist-State"
$ws.Cells.Item(39, 5).Value = 3
$ws.Cells.Item(39, 6).Value = 2.5
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(39, 8).Value = 3
$ws.Cells.Item(39, 9).Value = 3
$ws.Cells.Item(39, 10).Value = 1.5

$ws.Cells.Item(40, 4).Value = "#anarchist_communism"
$ws.Cells.Item(40, 5).Value = 1
$ws.Cells.Item(40, 6).Value = 2
$ws.Cells.Item(40, 7).Value = 4
$ws.Cells.Item(40, 8).Value = 3
$ws.Cells.Item(40, 9).Value = 3
$ws.Cells.Item(40, 10).Value = 1

$ws.Cells.Item(41, 4).Value = "#Conservative"
$ws.Cells.Item(41, 5).Value = 2
$ws.Cells.Item(41, 6).Value = 3
$ws.Cells.Item(41, 7).Value = 2
$ws.Cells.Item(41, 8).Value = 2
$ws.Cells.Item(41, 9).Value = 2.5
$ws.Cells.Item(41, 10).Value = 3

$ws.Cells.Item(42, 4).Value = "#Autocracy"
$ws.Cells.Item(42, 5).Value = 3
$ws.Cells.Item(42, 6).Value = 4
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = 1.5
$ws.Cells.Item(42, 9).Value = 1.5
$ws.Cells.Item(42, 10).Value = 3

$ws.Cells.Item(43, 4).Value = "#Mod_Vilayat_e_Faqih"
$ws.Cells.Item(43, 5).Value = 2
$ws.Cells.Item(43, 6).Value = 3
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(43, 8).Value = 2
$ws.Cells.Item(43, 9).Value = 2.5
$ws.Cells.Item(43, 10).Value = 2.5

$ws.Cells.Item(44, 4).Value = "#Vilayat_e_Faqih"
$ws.Cells.Item(44, 5).Value = 3.5
$ws.Cells.Item(44, 6).Value = 3
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(44, 8).Value = 2
$ws.Cells.Item(44, 9).Value = 2
$ws.Cells.Item(44, 10).Value = 2

$ws.Cells.Item(45, 4).Value = "#Kingdom"
$ws.Cells.Item(45, 5).Value = 2
$ws.Cells.Item(45, 6).Value = 4
$ws.Cells.Item(45, 7).Value = 2
$ws.Cells.Item(45, 8).Value = 2
$ws.Cells.Item(45, 9).Value = 2
$ws.Cells.Item(45, 10).Value = 3

$ws.Cells.Item(46, 4).Value = "#Caliphate"
$ws.Cells.Item(46, 5).Value = 2
$ws.Cells.Item(46, 6).Value = 4
$ws.Cells.Item(46, 7).Value = 1.5
$ws.Cells.Item(46, 8).Value = 1.5
$ws.Cells.Item(46, 9).Value = 1.5
$ws.Cells.Item(46, 10).Value = 4

$ws.Cells.Item(47, 4).Value = "#Neutral_Muslim_Brotherhood"
$ws.Cells.Item(47, 5).Value = 1.5
$ws.Cells.Item(47, 6).Value = 3
$ws.Cells.Item(47, 7).Value = 2.5
$ws.Cells.Item(47, 8).Value = 2.5
$ws.Cells.Item(47, 9).Value = 3
$ws.Cells.Item(47, 10).Value = 2.5

$ws.Cells.Item(48, 4).Value = "#Neutral_Autocracy"
$ws.Cells.Item(48, 5).Value = 3.5
$ws.Cells.Item(48, 6).Value = 3.5
$ws.Cells.Item(48, 7).Value = 2
$ws.Cells.Item(48, 8).Value = 2
$ws.Cells.Item(48, 9).Value = 2
$ws.Cells.Item(48, 10).Value = 2

$ws.Cells.Item(49, 4).Value = "#Neutral_conservatism"
$ws.Cells.Item(49, 5).Value = 2
$ws.Cells.Item(49, 6).Value = 3
$ws.Cells.Item(49, 7).Value = 2
$ws.Cells.Item(49, 8).Value = 2
$ws.Cells.Item(49, 9).Value = 2.5
$ws.Cells.Item(49, 10).Value = 2.5

$ws.Cells.Item(50, 4).Value = "#oligarchism"
$ws.Cells.Item(50, 5).Value = 3
$ws.Cells.Item(50, 6).Value = 3
$ws.Cells.Item(50, 7).Value = 2.5
$ws.Cells.Item(50, 8).Value = 2
$ws.Cells.Item(50, 9).Value = 2
$ws.Cells.Item(50, 10).Value = 2.5

$ws.Cells.Item(51, 4).Value = "#Neutral_Libertarian"
$ws.Cells.Item(51, 5).Value = 1.5
$ws.Cells.Item(51, 6).Value = 2
$ws.Cells.Item(51, 7).Value = 3.5
$ws.Cells.Item(51, 8).Value = 3
$ws.Cells.Item(51, 9).Value = 3
$ws.Cells.Item(51, 10).Value = 2

$ws.Cells.Item(52, 4).Value = "#Neutral_green"
$ws.Cells.Item(52, 5).Value = 2.5
$ws.Cells.Item(52, 6).Value = 1.5
$ws.Cells.Item(52, 7).Value = 4
$ws.Cells.Item(52, 8).Value = 2.5
$ws.Cells.Item(52, 9).Value = 2.5
$ws.Cells.Item(52, 10).Value = 2

$ws.Cells.Item(53, 4).Value = "#neutral_Social"
$ws.Cells.Item(53, 5).Value = 2.5
$ws.Cells.Item(53, 6).Value = 1.5
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 3
$ws.Cells.Item(53, 9).Value = 3
$ws.Cells.Item(53, 10).Value = 2

$ws.Cells.Item(54, 4).Value = "#Neutral_Communism"
$ws.Cells.Item(54, 5).Value = 3
$ws.Cells.Item(54, 6).Value = 2.5
$ws.Cells.Item(54, 7).Value = 2
$ws.Cells.Item(54, 8).Value = 2.5
$ws.Cells.Item(54, 9).Value = 2.5
$ws.Cells.Item(54, 10).Value = 2.5

$ws.Cells.Item(55, 4).Value = "#Nat_Populism"
$ws.Cells.Item(55, 5).Value = 1.5
$ws.Cells.Item(55, 6).Value = 4
$ws.Cells.Item(55, 7).Value = 1.5
$ws.Cells.Item(55, 8).Value = 2
$ws.Cells.Item(55, 9).Value = 2
$ws.Cells.Item(55, 10).Value = 4

$ws.Cells.Item(56, 4).Value = "#Nat_Fascism"
$ws.Cells.Item(56, 5).Value = 2.5
$ws.Cells.Item(56, 6).Value = 3.5
$ws.Cells.Item(56, 7).Value = 1.5
$ws.Cells.Item(56, 8).Value = 1.5
$ws.Cells.Item(56, 9).Value = 1.5
$ws.Cells.Item(56, 10).Value = 4

$ws.Cells.Item(57, 4).Value = "#Nat_Autocracy"
$ws.Cells.Item(57, 5).Value = 2.5
$ws.Cells.Item(57, 6).Value = 3
$ws.Cells.Item(57, 7).Value = 2
$ws.Cells.Item(57, 8).Value = 2
$ws.Cells.Item(57, 9).Value = 2
$ws.Cells.Item(57, 10).Value = 3.5

$ws.Cells.Item(58, 4).Value = "#Monarchist"
$ws.Cells.Item(58, 5).Value = 3
$ws.Cells.Item(58, 6).Value = 3
$ws.Cells.Item(58, 7).Value = 2
$ws.Cells.Item(58, 8).Value = 2
$ws.Cells.Item(58, 9).Value = 2
$ws.Cells.Item(58, 10).Value = 2

# Column K: row total, entered once then filled down (matches shared-formula pattern)
$ws.Range("K35").Formula = "=SUM(E35:J35)"
$ws.Range("K36:K58").Formula = "=SUM(E36:J36)"

# Update the visible selection to match the authored view state
$ws.Range("H35").Select()
